$wb = $excel.ActiveWorkbook

$iac = $wb.Worksheets.Item("IAC")

# Number of WTGs on IAC section increased from 66 to 132 (row 9, col E)
$iac.Range("E9").Value = 132

# Row 11: Nominal Current [A] on IAC Section
$iac.Range("E11").Value = 92.08
$iac.Range("F11").Formula = '=F10*$E$11'
$iac.Range("G11").Formula = '=G10*$E$11'
$iac.Range("H11").Formula = '=H10*$E$11'

# Row 12: WTG Nominal Power [MVA]
$iac.Range("E12").Value = 1
$iac.Range("F12").Value = 1
$iac.Range("G12").Value = 1
$iac.Range("H12").Value = 1

# Row 13: Losses (% WTG Energy Yield/year)
$iac.Range("E13").Value = 0.101
$iac.Range("F13").Value = 0.101
$iac.Range("G13").Value = 0.101
$iac.Range("H13").Value = 0.101

# Column D (rows 18-46): IAC String Losses (MWh/year) = C / 0.95
$iac.Range("D18").Formula = "=C18/0.95"
$iac.Range("D19:D45").Formula = "=C19/0.95"
$iac.Range("D46").Formula = "=C46/0.95"

# Summary rows at bottom
$iac.Range("K49").Formula = "='WTG Yield Wake'!J43*4"
$iac.Range("K50").Formula = "=(K48/K49)"

$wb.Application.Calculate()
